$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (marker G2=5489)
$ws.Range("H2").Value = 840.8333
$ws.Range("I2").Value = 817
$ws.Range("J2").Value = 864.6667
$ws.Range("K2").Value = 817
$ws.Range("L2").Value = 864.6667
$ws.Range("M2").Value = -704
$ws.Range("N2").Value = -1090.6667
# Row 9 (marker G9=5487)
$ws.Range("H9").Value = 274.5
$ws.Range("I9").Value = 162.5
$ws.Range("K9").Value = 162.5
$ws.Range("M9").Value = 6.5
# Row 10 (marker G10=1959)
$ws.Range("H10").Value = 26232.6
$ws.Range("I10").Value = 9999
$ws.Range("K10").Value = 9999
$ws.Range("M10").Value = -9706
# Row 18 (marker G18=5471)
$ws.Range("H18").Value = 375
$ws.Range("I18").Value = 350
$ws.Range("K18").Value = 350
$ws.Range("M18").Value = -66
# Row 75 (marker G75=10654)
$ws.Range("H75").Value = 38104.668
$ws.Range("J75").Value = 38104.668
$ws.Range("L75").Value = 38104.668
$ws.Range("N75").Value = -39976.668
# Row 78 (marker G78=10654)
$ws.Range("H78").Value = 38104.668
$ws.Range("J78").Value = 38104.668
$ws.Range("L78").Value = 114314.004
$ws.Range("N78").Value = -123674.004
# Row 111 (marker G111=27768)
$ws.Range("H111").Value = 1440.375
$ws.Range("I111").Value = 733.8570999999999
$ws.Range("J111").Value = 1989.8889
$ws.Range("K111").Value = 2201.5713
$ws.Range("L111").Value = 5969.6667
$ws.Range("M111").Value = 865.4287000000004
$ws.Range("N111").Value = -12103.6667
# Row 112 (marker G112=27960)
$ws.Range("H112").Value = 682133.75
$ws.Range("I112").Value = 2722
$ws.Range("K112").Value = 8166
$ws.Range("M112").Value = -7058
# Row 127 (marker G127=36114)
$ws.Range("H127").Value = 3476.8
$ws.Range("I127").Value = 1398.6
$ws.Range("J127").Value = 5555
$ws.Range("K127").Value = 4195.799999999999
$ws.Range("L127").Value = 16665
$ws.Range("M127").Value = 764.2000000000007
$ws.Range("N127").Value = -26585
# Row 129 (marker G129=36115)
$ws.Range("H129").Value = 1200
$ws.Range("I129").Value = 937.3077
$ws.Range("J129").Value = 1626.875
$ws.Range("K129").Value = 2811.9231
$ws.Range("L129").Value = 4880.625
$ws.Range("M129").Value = 2188.0769
$ws.Range("N129").Value = -14880.625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (marker G32=44147)
$ws.Range("H32").Value = 2905.3704
$ws.Range("I32").Value = 2905.3704
$ws.Range("K32").Value = 2905.3704
$ws.Range("M32").Value = -2618.3704
# Row 34 (marker G34=2753)
$ws.Range("H34").Value = 45304
$ws.Range("I34").Value = 45304
$ws.Range("K34").Value = 45304
$ws.Range("M34").Value = -45033
# Row 122 (marker G122=36168)
$ws.Range("H122").Value = 1305.375
$ws.Range("I122").Value = 1175.0834
$ws.Range("K122").Value = 3525.2502
$ws.Range("M122").Value = -1075.2502
# Row 131 (marker G131=34706)
$ws.Range("H131").Value = 88000
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
# Row 132 (marker G132=43997)
$ws.Range("H132").Value = 2582.6785
$ws.Range("I132").Value = 2568.68
$ws.Range("K132").Value = 7706.039999999999
$ws.Range("M132").Value = -5176.039999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (marker G3=27713)
$ws.Range("H3").Value = 1921.7826
$ws.Range("I3").Value = 1533.0588
$ws.Range("K3").Value = 1533.0588
$ws.Range("M3").Value = -1419.0588
# Row 14 (marker G14=2342)
$ws.Range("H14").Value = 699
$ws.Range("I14").Value = 699
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 699
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -527
$ws.Range("N14").ClearContents()
# Row 103 (marker G103=18514)
$ws.Range("H103").Value = 23176.25
$ws.Range("J103").Value = 23176.25
$ws.Range("L103").Value = 23176.25
$ws.Range("N103").Value = -25520.25
# Row 128 (marker G128=38749)
$ws.Range("H128").Value = 8040
$ws.Range("I128").Value = 8040
$ws.Range("K128").Value = 24120
$ws.Range("M128").Value = -21630
# Row 134 (marker G134=43998)
$ws.Range("H134").Value = 3012.818
$ws.Range("I134").Value = 3012.818
$ws.Range("K134").Value = 9038.454000000002
$ws.Range("M134").Value = -6503.454000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22 (marker G22=5367)
$ws.Range("H22").Value = 613.2
$ws.Range("I22").Value = 346.75
$ws.Range("K22").Value = 346.75
$ws.Range("M22").Value = 3.25
# Row 55 (marker G55=1855)
$ws.Range("H55").Value = 25745.2
$ws.Range("I55").Value = 24613.143
$ws.Range("K55").Value = 24613.143
$ws.Range("M55").Value = -24298.143
# Row 99 (marker G99=36198)
$ws.Range("H99").Value = 1686.7142
$ws.Range("I99").Value = 1461.091
$ws.Range("K99").Value = 1461.091
$ws.Range("M99").Value = 36.90900000000011
# Row 126 (marker G126=36198)
$ws.Range("H126").Value = 1686.7142
$ws.Range("I126").Value = 1461.091
$ws.Range("K126").Value = 4383.272999999999
$ws.Range("M126").Value = -1913.272999999999
# Row 132 (marker G132=44019)
$ws.Range("H132").Value = 3329.9285
$ws.Range("I132").Value = 3427.5
$ws.Range("K132").Value = 10282.5
$ws.Range("M132").Value = -7752.5
# Row 134 (marker G134=44020)
$ws.Range("H134").Value = 1657.95
$ws.Range("I134").Value = 1376.8422
$ws.Range("K134").Value = 4130.5266
$ws.Range("M134").Value = -1595.5266

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 37 (marker G37=9516)
$ws.Range("H37").Value = 105653.71
$ws.Range("J37").Value = 105653.71
$ws.Range("L37").Value = 316961.13
$ws.Range("N37").Value = -317185.13
# Row 129 (marker G129=36054)
$ws.Range("H129").Value = 1741.5
$ws.Range("I129").Value = 495.625
$ws.Range("J129").Value = 2572.0833
$ws.Range("K129").Value = 1486.875
$ws.Range("L129").Value = 7716.249899999999
$ws.Range("M129").Value = 3513.125
$ws.Range("N129").Value = -17716.2499

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 15 (marker G15=12018)
$ws.Range("H15").Value = 9827.5
$ws.Range("J15").Value = 9827.5
$ws.Range("L15").Value = 9827.5
$ws.Range("N15").Value = -10403.5
# Row 26 (marker G26=4254)
$ws.Range("H26").Value = 13376.4
$ws.Range("I26").Value = 7899
$ws.Range("K26").Value = 7899
$ws.Range("M26").Value = -7619
# Row 48 (marker G48=4337)
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 50 (marker G50=4254)
$ws.Range("H50").Value = 13376.4
$ws.Range("I50").Value = 7899
$ws.Range("K50").Value = 7899
$ws.Range("M50").Value = -7401
# Row 80 (marker G80=12521)
$ws.Range("H80").Value = 16557.334
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 16557.334
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 16557.334
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -18553.334
# Row 81 (marker G81=12018)
$ws.Range("H81").Value = 9827.5
$ws.Range("J81").Value = 9827.5
$ws.Range("L81").Value = 9827.5
$ws.Range("N81").Value = -11823.5
# Row 83 (marker G83=12521)
$ws.Range("H83").Value = 16557.334
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 16557.334
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 82786.67
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -92770.67
# Row 84 (marker G84=12018)
$ws.Range("H84").Value = 9827.5
$ws.Range("J84").Value = 9827.5
$ws.Range("L84").Value = 29482.5
$ws.Range("N84").Value = -39466.5
# Row 95 (marker G95=18235)
$ws.Range("H95").Value = 28114.334
$ws.Range("J95").Value = 28114.334
$ws.Range("L95").Value = 28114.334
$ws.Range("N95").Value = -33606.334
# Row 122 (marker G122=36182)
$ws.Range("H122").Value = 7246.4
# Row 126 (marker G126=36184)
$ws.Range("H126").Value = 12149.7
$ws.Range("I126").Value = 9749.25
$ws.Range("K126").Value = 29247.75
$ws.Range("M126").Value = -26777.75
# Row 132 (marker G132=44008)
$ws.Range("H132").Value = 4019.8
$ws.Range("I132").Value = 3794.3076
$ws.Range("K132").Value = 11382.9228
$ws.Range("M132").Value = -8852.9228

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 43 (marker G43=4314)
$ws.Range("H43").Value = 38099.8
$ws.Range("J43").Value = 17624.75
$ws.Range("L43").Value = 17624.75
$ws.Range("N43").Value = -18010.75
# Row 55 (marker G55=5284)
$ws.Range("H55").Value = 315.48276
$ws.Range("I55").Value = 334
$ws.Range("J55").Value = 309.5909
$ws.Range("K55").Value = 334
$ws.Range("L55").Value = 309.5909
$ws.Range("M55").Value = -161
$ws.Range("N55").Value = -655.5908999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 25 (marker G25=3064)
$ws.Range("H25").Value = 30000
$ws.Range("J25").Value = 30000
$ws.Range("L25").Value = 30000
$ws.Range("N25").Value = -30586
# Row 40 (marker G40=3601)
$ws.Range("H40").Value = 14998.667
$ws.Range("J40").Value = 14998.667
$ws.Range("L40").Value = 14998.667
$ws.Range("N40").Value = -15296.667
# Row 41 (marker G41=21725)
$ws.Range("H41").Value = 193657.75
$ws.Range("I41").Value = 48687
$ws.Range("J41").Value = 241981.33
$ws.Range("K41").Value = 48687
$ws.Range("L41").Value = 241981.33
$ws.Range("M41").Value = -48297
$ws.Range("N41").Value = -242761.33
# Row 81 (marker G81=12596)
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
# Row 84 (marker G84=12596)
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
# Row 122 (marker G122=36208)
$ws.Range("H122").Value = 2173.2144
$ws.Range("J122").Value = 2002
$ws.Range("L122").Value = 6006
$ws.Range("N122").Value = -10906
# Row 126 (marker G126=36210)
$ws.Range("H126").Value = 1234.7778
$ws.Range("I126").Value = 1234.7778
$ws.Range("K126").Value = 3704.3334
$ws.Range("M126").Value = -1234.3334
# Row 130 (marker G130=34705)
$ws.Range("H130").Value = 44443
$ws.Range("J130").Value = 44443
$ws.Range("L130").Value = 44443
$ws.Range("N130").Value = -54483
# Row 135 (marker G135=42043)
$ws.Range("H135").Value = 55167.25
$ws.Range("J135").Value = 55167.25
$ws.Range("L135").Value = 55167.25
$ws.Range("N135").Value = -65307.25
# Row 136 (marker G136=44031)
$ws.Range("H136").Value = 695.2143
$ws.Range("I136").Value = 695.2143
$ws.Range("K136").Value = 2085.6429
$ws.Range("M136").Value = 464.3571000000002
